$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 594.5333000000001
$ws.Range("I15").Value = 594.5333000000001
$ws.Range("K15").Value = 1783.5999
$ws.Range("M15").Value = -1614.5999

$ws.Range("H137").Value = 3466.25
$ws.Range("I137").Value = 3286
$ws.Range("K137").Value = 9858
$ws.Range("M137").Value = -7308

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()

$ws.Range("H11").Value = 2500
$ws.Range("J11").Value = 2500
$ws.Range("L11").Value = 2500
$ws.Range("N11").Value = -2788

$ws.Range("H13").Value = 14398.6
$ws.Range("J13").Value = 14398.6
$ws.Range("L13").Value = 14398.6
$ws.Range("N13").Value = -14686.6

$ws.Range("H32").Value = 3412.84
$ws.Range("I32").Value = 3057.4348
$ws.Range("K32").Value = 3057.4348
$ws.Range("M32").Value = -2770.4348

$ws.Range("H37").Value = 18330.555
$ws.Range("J37").Value = 19996.875
$ws.Range("L37").Value = 19996.875
$ws.Range("N37").Value = -20542.875

$ws.Range("H55").Value = 24998.75
$ws.Range("J55").Value = 24998.75
$ws.Range("L55").Value = 24998.75
$ws.Range("N55").Value = -25628.75

$ws.Range("H74").Value = 5000
$ws.Range("I74").Value = 5000
$ws.Range("K74").Value = 5000
$ws.Range("M74").Value = -4126

$ws.Range("H77").Value = 5000
$ws.Range("I77").Value = 5000
$ws.Range("K77").Value = 25000
$ws.Range("M77").Value = -20632

$ws.Range("H80").Value = 37499.375
$ws.Range("J80").Value = 37142.145
$ws.Range("L80").Value = 37142.145
$ws.Range("N80").Value = -39138.145

$ws.Range("H83").Value = 37499.375
$ws.Range("J83").Value = 37142.145
$ws.Range("L83").Value = 111426.435
$ws.Range("N83").Value = -121410.435

$ws.Range("H110").Value = 839.7143
$ws.Range("I110").Value = 813
$ws.Range("J110").Value = 1000
$ws.Range("K110").Value = 813
$ws.Range("L110").Value = 1000
$ws.Range("M110").Value = 1232
$ws.Range("N110").Value = -5090

$ws.Range("H122").Value = 4148
$ws.Range("I122").Value = 4148
$ws.Range("K122").Value = 12444
$ws.Range("M122").Value = -9994

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 14970
$ws.Range("J35").Value = 14970
$ws.Range("L35").Value = 14970
$ws.Range("N35").Value = -15590

$ws.Range("H82").Value = 30621.666
$ws.Range("I82").Value = 16749.285
$ws.Range("J82").Value = 39449.547
$ws.Range("K82").Value = 16749.285
$ws.Range("L82").Value = 39449.547
$ws.Range("M82").Value = -16366.285
$ws.Range("N82").Value = -40215.547

$ws.Range("H85").Value = 30621.666
$ws.Range("I85").Value = 16749.285
$ws.Range("J85").Value = 39449.547
$ws.Range("K85").Value = 16749.285
$ws.Range("L85").Value = 39449.547
$ws.Range("M85").Value = -15423.285
$ws.Range("N85").Value = -42101.547

$ws.Range("H86").Value = 4925.5
$ws.Range("I86").Value = 4352.273
$ws.Range("J86").Value = 5826.2856
$ws.Range("K86").Value = 4352.273
$ws.Range("L86").Value = 5826.2856
$ws.Range("M86").Value = -3229.273
$ws.Range("N86").Value = -8072.2856

$ws.Range("H89").Value = 4925.5
$ws.Range("I89").Value = 4352.273
$ws.Range("J89").Value = 5826.2856
$ws.Range("K89").Value = 21761.365
$ws.Range("L89").Value = 29131.428
$ws.Range("M89").Value = -16145.365
$ws.Range("N89").Value = -40363.428

$ws.Range("H94").Value = 1489.1818
$ws.Range("I94").Value = 1519.3334
$ws.Range("K94").Value = 1519.3334
$ws.Range("M94").Value = -1068.3334

$ws.Range("H105").Value = 4549.4287
$ws.Range("J105").Value = 3074
$ws.Range("L105").Value = 3074
$ws.Range("N105").Value = -6568

$ws.Range("H128").Value = 10000
$ws.Range("I128").Value = 10000
$ws.Range("K128").Value = 30000
$ws.Range("M128").Value = -27510

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 19998.75
$ws.Range("J41").Value = 19998.75
$ws.Range("L41").Value = 19998.75
$ws.Range("N41").Value = -20854.75

$ws.Range("H50").Value = 28425.428
$ws.Range("J50").Value = 29649.166
$ws.Range("L50").Value = 29649.166
$ws.Range("N50").Value = -30899.166

$ws.Range("H59").Value = 33798.332
$ws.Range("J59").Value = 33926.785
$ws.Range("L59").Value = 33926.785
$ws.Range("N59").Value = -36216.785

$ws.Range("H60").Value = 22512.572
$ws.Range("J60").Value = 24582.5
$ws.Range("L60").Value = 24582.5
$ws.Range("N60").Value = -25604.5

$ws.Range("H68").Value = 37807.777
$ws.Range("J68").Value = 40000
$ws.Range("L68").Value = 40000
$ws.Range("N68").Value = -41498

$ws.Range("H71").Value = 37807.777
$ws.Range("J71").Value = 40000
$ws.Range("L71").Value = 120000
$ws.Range("N71").Value = -127488

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H134").Value = 3115
$ws.Range("I134").Value = 2852.7144
$ws.Range("K134").Value = 8558.143199999999
$ws.Range("M134").Value = -6023.143199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 2000
$ws.Range("J52").Value = 2000
$ws.Range("L52").Value = 6000
$ws.Range("N52").Value = -6532

$ws.Range("H60").Value = 339.33334
$ws.Range("I60").Value = 339.33334
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 1018.00002
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -767.0000200000001
$ws.Range("N60").ClearContents()

$ws.Range("H81").Value = 5500
$ws.Range("J81").Value = 5500
$ws.Range("L81").Value = 16500
$ws.Range("N81").Value = -18746

$ws.Range("H84").Value = 5500
$ws.Range("J84").Value = 5500
$ws.Range("L84").Value = 49500
$ws.Range("N84").Value = -60732

$ws.Range("H109").Value = 1618.75
$ws.Range("I109").Value = 481.66666
$ws.Range("K109").Value = 1444.99998
$ws.Range("M109").Value = -404.9999800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 11875
$ws.Range("I43").Value = 2500
$ws.Range("J43").Value = 15000
$ws.Range("K43").Value = 2500
$ws.Range("L43").Value = 15000
$ws.Range("M43").Value = -2349
$ws.Range("N43").Value = -15302

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 2999.5
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 1006
$ws.Range("J10").Value = 1006
$ws.Range("L10").Value = 1006
$ws.Range("N10").Value = -1344

$ws.Range("H132").Value = 1347
$ws.Range("I132").Value = 1364.3334
$ws.Range("K132").Value = 4093.0002
$ws.Range("M132").Value = -1563.0002
